$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update the last row (2025Q3) with the refreshed recurrence metrics
$ws.Range("C29").Value = 156
$ws.Range("D29").Value = 24
$ws.Range("E29").Value = 132
$ws.Range("F29").Value = 4.130808950086059
